$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values in rows 2-5
$ws.Range("C2").Value = -2.717172

$ws.Range("B3").Value = 0.000735
$ws.Range("C3").Value = -2.668100

$ws.Range("B4").Value = 0.005020
$ws.Range("C4").Value = -2.815317
$ws.Range("D4").Value = 0.147217

$ws.Range("B5").Value = 0.009305
$ws.Range("C5").Value = -3.698622
$ws.Range("D5").Value = 1.030522

# Remove rows 6-8 entirely (data now ends at row 5)
$ws.Range("A6:D8").Delete()
